$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.023.86"
$ws.Range("E2").Value = "  -1.42%  "
$ws.Range("D3").Value = "1.829.67"
$ws.Range("E3").Value = "  -1.39%  "
$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'238.92"
$ws.Range("E5").Value = "  -2.35%  "
$ws.Range("D6").Value = "'0.6659"
$ws.Range("E6").Value = "  -3.91%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "'0.2950"
$ws.Range("E8").Value = "  -3.41%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.07383"
$ws.Range("E9").Value = "  -3.52%  "
$ws.Range("D10").Value = "'22.69"
$ws.Range("E10").Value = "  -3.54%  "
$ws.Range("D11").Value = "'0.07637"
$ws.Range("E11").Value = "  -1.59%  "
$ws.Range("D12").Value = "1.839.72"
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("D13").Value = "'5.007"
$ws.Range("E13").Value = "  -2.34%  "
$ws.Range("D14").Value = "'0.6717"
$ws.Range("E14").Value = "  -2.57%  "
$ws.Range("D15").Value = "'86.09"
$ws.Range("E15").Value = "  -4.76%  "
$ws.Range("D16").Value = "'6.114"
$ws.Range("E16").Value = "  -4.82%  "
$ws.Range("D17").Value = "29.039.01"
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("D18").Value = "'0.000008205"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").Value = "'227.19"
$ws.Range("E19").Value = "  -4.25%  "
$ws.Range("D20").Value = "'12.43"
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("D21").Value = "'0.9995"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").Value = "'7.268"
$ws.Range("E22").Value = "  -4.63%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "'160.35"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("D25").Value = "'0.1422"
$ws.Range("E25").Value = "  -4.50%  "
$ws.Range("D26").Value = "'8.654"
$ws.Range("D27").Value = "'17.93"
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("D28").Value = "'1.502"
$ws.Range("E28").Value = "  -2.04%  "
$ws.Range("D29").Value = "'4.224"
$ws.Range("E29").Value = "  -0.58%  "
$ws.Range("D30").Value = "'4.112"
$ws.Range("E30").Value = "  -0.83%  "
$ws.Range("D31").Value = "'1.196"
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("D32").Value = "'0.05363"
$ws.Range("E32").Value = "  +5.04%  "
$ws.Range("D33").Value = "'0.7470"
$ws.Range("E33").Value = "  -2.60%  "
$ws.Range("D34").Value = "'1.843"
$ws.Range("E34").Value = "  -2.12%  "
$ws.Range("D35").Value = "'1.125"
$ws.Range("E35").Value = "  -2.05%  "
$ws.Range("D36").Value = "'2.679"
$ws.Range("E36").Value = "  -0.31%  "
$ws.Range("D37").Value = "1.294.29"
$ws.Range("E37").Value = "  -2.48%  "
$ws.Range("D38").Value = "'0.01804"
$ws.Range("E38").Value = "  -2.94%  "
$ws.Range("D39").Value = "'2.703"
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("D40").Value = "'0.9209"
$ws.Range("E40").Value = "  -6.00%  "
$ws.Range("D41").Value = "'6.016"
$ws.Range("E41").Value = "  +3.26%  "
$ws.Range("D42").Value = "'103.89"
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("D43").Value = "'0.9986"
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("D44").Value = "'0.07867"
$ws.Range("E44").Value = "  +18.94%  "
$ws.Range("D45").Value = "1.978.66"
$ws.Range("E45").Value = "  -1.24%  "
$ws.Range("E46").Value = "  -3.48%  "
$ws.Range("D47").Value = "'0.5176"
$ws.Range("E47").Value = "  -0.84%  "
$ws.Range("D48").Value = "'63.32"
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("D49").Value = "'1.748"
$ws.Range("E49").Value = "  -1.29%  "
$ws.Range("D50").Value = "'9.258"
$ws.Range("E50").Value = "  -5.27%  "
$ws.Range("D51").Value = "'0.05908"
$ws.Range("E51").Value = "  -0.22%  "
